$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing-address paragraph "989 Story Road, San Jose CA 95122"
#    (the one in the free-flowing letter body, not the one inside the details
#    table) into two paragraphs: "989 Story Road" and "San Jose, CA 95122".
#    Locate it precisely via the Paragraphs collection so the table's copy of
#    the same text is left untouched.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "989 Story Road, San Jose CA 95122`r" -and $para.Range.Information(12) -eq $false) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne -1) {
    $addrPara = $d.Paragraphs.Item($targetIndex)
    $addrRange = $addrPara.Range
    $addrRange.Find.Execute("989 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false, $true, 1, $false, "989 Story Road`rSan Jose, CA 95122", 2) | Out-Null

    # Restore run formatting (Arial 11pt, incl. complex-script sizing) on the
    # newly created second line, matching the sibling address paragraph.
    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newRange = $newPara.Range
    $newRange.Font.Name = "Arial"
    $newRange.Font.NameBi = "Arial"
    $newRange.Font.Size = 11
    $newRange.Font.SizeBi = 11
}

# 3. Remove the empty "No Spacing" paragraph directly after the
#    "...Board of Directors" signature line.
$bodIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Vietnam Town Condominium Owners Association Board of Directors`r") {
        $bodIndex = $i
        break
    }
}

if ($bodIndex -ne -1) {
    $emptyPara = $d.Paragraphs.Item($bodIndex + 1)
    if ($emptyPara.Range.Text -eq "`r") {
        $emptyPara.Range.Delete()
    }
}
